# Update the NATMI Thbs1-Sdc1 ligand-receptor table with newly recomputed TPM-based values.
#
# Column layout (row 1 headers):
#  A Sending cluster          B Ligand symbol            C Receptor symbol        D Target cluster
#  E Ligand-expressing cells  F Ligand detection rate
#  G Ligand average expression value        H Ligand total expression value
#  I Ligand derived specificity (avg)       J Ligand derived specificity (total)
#  K Receptor-expressing cells              L Receptor detection rate
#  M Receptor average expression value      N Receptor total expression value
#  O Receptor derived specificity (avg)     P Receptor derived specificity (total)
#  Q Edge average expression weight         R Edge total expression weight
#  S Edge average expression derived specificity   T Edge total expression derived specificity
#
# G,H,I,J depend only on the Sending cluster (column A); K,L,M,N,O,P depend only on
# the Target cluster (column D); Q=G*M, R=H*N, S=I*O, T=J*P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ligand-side values (Sending cluster -> G, H, I, J)
$sendingMap = @{
    "ECs"               = @(16.71895933333333, 50.156878,        0.02912144738161902, 0.03059269312988411)
    "FAPs"              = @(155.500389,        466.501167,       0.2708539632042961,  0.2845377865576845)
    "Inflammatory-Mac"  = @(194.8548433333333, 584.56453,        0.3394024086099587,  0.3565493705749576)
    "MuSCs"             = @(82.82950199999999, 165.659004,       0.1442742299952585,  0.1010420758958371)
    "Resolving-Mac"     = @(124.2078576666667, 372.623573,       0.2163479508088675,  0.2272780738416368)
}

# New receptor-side values (Target cluster -> K, L, M, N, O, P)
$targetMap = @{
    "ECs"               = @(3, 1,                  1.306376666666667,  3.91913,   0.06159635513812315, 0.07271399171915481)
    "FAPs"              = @(3, 1,                  8.648731,           25.946193, 0.4077922698431246,  0.4813954277979023)
    "Inflammatory-Mac"  = @(2, 0.6666666666666666, 0.637617,           1.912851,  0.03006398091472189, 0.03549028273468269)
    "MuSCs"             = @(2, 1,                  9.728125,           19.45625,  0.4586862714388558,  0.3609835859963323)
    "Resolving-Mac"     = @(3, 1,                  0.8878186666666666, 2.663456,  0.0418611226651744,  0.0494167117519279)
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value2
    $target  = $ws.Cells.Item($r, 4).Value2

    $sv = $sendingMap[$sending]
    $tv = $targetMap[$target]

    $g = $sv[0]; $h = $sv[1]; $i = $sv[2]; $j = $sv[3]
    $k = $tv[0]; $l = $tv[1]; $m = $tv[2]; $n = $tv[3]; $o = $tv[4]; $p = $tv[5]

    $ws.Cells.Item($r, 7).Value  = $g   # G Ligand average expression value
    $ws.Cells.Item($r, 8).Value  = $h   # H Ligand total expression value
    $ws.Cells.Item($r, 9).Value  = $i   # I Ligand derived specificity (avg)
    $ws.Cells.Item($r, 10).Value = $j   # J Ligand derived specificity (total)

    $ws.Cells.Item($r, 11).Value = $k   # K Receptor-expressing cells
    $ws.Cells.Item($r, 12).Value = $l   # L Receptor detection rate
    $ws.Cells.Item($r, 13).Value = $m   # M Receptor average expression value
    $ws.Cells.Item($r, 14).Value = $n   # N Receptor total expression value
    $ws.Cells.Item($r, 15).Value = $o   # O Receptor derived specificity (avg)
    $ws.Cells.Item($r, 16).Value = $p   # P Receptor derived specificity (total)

    $ws.Cells.Item($r, 17).Value = $g * $m   # Q Edge average expression weight
    $ws.Cells.Item($r, 18).Value = $h * $n   # R Edge total expression weight
    $ws.Cells.Item($r, 19).Value = $i * $o   # S Edge average expression derived specificity
    $ws.Cells.Item($r, 20).Value = $j * $p   # T Edge total expression derived specificity
}
